# Size Estimating Template.xlsx edit script
# Applies the changes described by the commit "Subí PSP Clase Usuario y BD_Usuario":
#   - sheet renamed from excel(1) -> excel(6) (and its local defined name excel_1 -> excel_6)
#   - the web-query host port moved from 2470 to 2469 (reflected in the hyperlink display text)
#   - "PROBE estimating basis used" cell changed from D to C
#   - a few size values bumped from 21 to 27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet and keep the local defined name ("excel_<n>") in sync with it.
$ws.Name = "excel(6)"
$defName = $wb.Names.Item(1)
$defName.Name = "excel_6"

# 2) Update the cells that hold numeric size values (21 -> 27).
$ws.Range("J16").Value = 27
$ws.Range("J17").Value = 27
$ws.Range("C24").Value = 27

# 3) "PROBE estimating basis used: (A, B, C, or D)" answer cell changes from D to C.
$ws.Range("D30").Value = "C"
$ws.Range("F30").Value = "C"

# 4) The reporting host port changed from 2470 to 2469; the hyperlink display text
#    (and the underlying connection/query URLs) reflect the new port.
foreach ($hl in $ws.Hyperlinks) {
    $hl.TextToDisplay = $hl.Address -replace "localhost:2470", "localhost:2469"
}
